$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update VillagerA's Speed value (G4) from 2.0 to 1.5
$ws.Range("G4").Value = 1.5

# Add new row 5 for VillagerB
$ws.Range("A5").Value = "VillagerB"
$ws.Range("B5").Value = 100.0
$ws.Range("C5").Value = "Assets/Prefabs/NPC/VillagerB.prefab"
$ws.Range("D5").Value = "VillagerStats"
$ws.Range("E5").Value = 40.0
$ws.Range("F5").Value = 40.0
$ws.Range("G5").Value = 1.0
$ws.Range("H5").Value = 1.0
$ws.Range("I5").Value = "VillagerB"
